# Cost Estimate.xlsx - BOM sheet update
#
# - Corrected order quantities for three line items (column E):
#     row 3 (C1005C0G1H100D050BA):  4 -> 1
#     row 4 (C1005X5R1A105K050BB):  1 -> 6
#     row 6 (ERJ-2GE0R00X):         3 -> 4
# - Highlighted the "Description" cells of the rows that were touched
#   (and a few related rows) with a solid yellow fill so the change
#   stands out for review.
# - Left the cursor sitting on B8 (the row being double-checked).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Quantity corrections -------------------------------------------------
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 6
$ws.Range("E6").Value = 4

# --- Highlight the affected part descriptions -----------------------------
$ws.Range("D2").Interior.Color = 65535
$ws.Range("D3").Interior.Color = 65535
$ws.Range("D4").Interior.Color = 65535
$ws.Range("D6").Interior.Color = 65535
$ws.Range("D7").Interior.Color = 65535
$ws.Range("D9").Interior.Color = 65535
$ws.Range("D10").Interior.Color = 65535
$ws.Range("D11").Interior.Color = 65535

# --- Leave the selection on B8 ---------------------------------------------
[void]$ws.Range("B8").Select()
